$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - PACOTE PRE-OPERATORIO PEDIATRICO OTORRINO
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 0

# Row 5 - ADENOIDECTOMIA PEDIATRICO
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("K5").Value = 0

# Row 6 - AMIGDALECTOMIA- PEDIATRICO
$ws.Range("D6").Value = 0

# Row 7 - AMIGDALECTOMIA COM ADENOIDECTOMIA - PEDIATRICO
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

# Row 17 - TOTAL
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
